$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.08661466666666667
$ws.Cells.Item(2, 8).Value = 0.259844
$ws.Cells.Item(2, 9).Value = 0.5374496355558498
$ws.Cells.Item(2, 10).Value = 0.5374496355558498
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.572065666666667
$ws.Cells.Item(2, 14).Value = 4.716197
$ws.Cells.Item(2, 15).Value = 0.1759712293834306
$ws.Cells.Item(2, 16).Value = 0.1759712293834305
$ws.Cells.Item(2, 17).Value = 0.1361639436964445
$ws.Cells.Item(2, 18).Value = 1.225475493268
$ws.Cells.Item(2, 19).Value = 0.09457567310043961
$ws.Cells.Item(2, 20).Value = 0.09457567310043959
$ws.Cells.Item(3, 7).Value = 0.08661466666666667
$ws.Cells.Item(3, 8).Value = 0.259844
$ws.Cells.Item(3, 9).Value = 0.5374496355558498
$ws.Cells.Item(3, 10).Value = 0.5374496355558498
$ws.Cells.Item(3, 15).Value = 0.4743638053196239
$ws.Cells.Item(3, 16).Value = 0.4743638053196239
$ws.Cells.Item(3, 17).Value = 0.3670557210146667
$ws.Cells.Item(3, 18).Value = 3.303501489132
$ws.Cells.Item(3, 19).Value = 0.254946654289918
$ws.Cells.Item(3, 20).Value = 0.254946654289918
$ws.Cells.Item(4, 7).Value = 0.08661466666666667
$ws.Cells.Item(4, 8).Value = 0.259844
$ws.Cells.Item(4, 9).Value = 0.5374496355558498
$ws.Cells.Item(4, 10).Value = 0.5374496355558498
$ws.Cells.Item(4, 13).Value = 3.123785
$ws.Cells.Item(4, 14).Value = 9.371354999999999
$ws.Cells.Item(4, 15).Value = 0.3496649652969456
$ws.Cells.Item(4, 16).Value = 0.3496649652969455
$ws.Cells.Item(4, 17).Value = 0.2705655965133333
$ws.Cells.Item(4, 18).Value = 2.43509036862
$ws.Cells.Item(4, 19).Value = 0.1879273081654923
$ws.Cells.Item(4, 20).Value = 0.1879273081654922
$ws.Cells.Item(5, 9).Value = 0.2407068810034004
$ws.Cells.Item(5, 10).Value = 0.2407068810034003
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.572065666666667
$ws.Cells.Item(5, 14).Value = 4.716197
$ws.Cells.Item(5, 15).Value = 0.1759712293834306
$ws.Cells.Item(5, 16).Value = 0.1759712293834305
$ws.Cells.Item(5, 17).Value = 0.06098357134133334
$ws.Cells.Item(5, 18).Value = 0.5488521420719999
$ws.Cells.Item(5, 19).Value = 0.04235748577121949
$ws.Cells.Item(5, 20).Value = 0.04235748577121948
$ws.Cells.Item(6, 9).Value = 0.2407068810034004
$ws.Cells.Item(6, 10).Value = 0.2407068810034003
$ws.Cells.Item(6, 15).Value = 0.4743638053196239
$ws.Cells.Item(6, 16).Value = 0.4743638053196239
$ws.Cells.Item(6, 19).Value = 0.1141826320393909
$ws.Cells.Item(6, 20).Value = 0.1141826320393909
$ws.Cells.Item(7, 9).Value = 0.2407068810034004
$ws.Cells.Item(7, 10).Value = 0.2407068810034003
$ws.Cells.Item(7, 13).Value = 3.123785
$ws.Cells.Item(7, 14).Value = 9.371354999999999
$ws.Cells.Item(7, 15).Value = 0.3496649652969456
$ws.Cells.Item(7, 16).Value = 0.3496649652969455
$ws.Cells.Item(7, 17).Value = 0.12117786772
$ws.Cells.Item(7, 18).Value = 1.09060080948
$ws.Cells.Item(7, 19).Value = 0.08416676319279
$ws.Cells.Item(7, 20).Value = 0.08416676319278997
$ws.Cells.Item(8, 7).Value = 0.035752
$ws.Cells.Item(8, 8).Value = 0.107256
$ws.Cells.Item(8, 9).Value = 0.2218434834407499
$ws.Cells.Item(8, 10).Value = 0.2218434834407499
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.572065666666667
$ws.Cells.Item(8, 14).Value = 4.716197
$ws.Cells.Item(8, 15).Value = 0.1759712293834306
$ws.Cells.Item(8, 16).Value = 0.1759712293834305
$ws.Cells.Item(8, 17).Value = 0.05620449171466667
$ws.Cells.Item(8, 18).Value = 0.505840425432
$ws.Cells.Item(8, 19).Value = 0.03903807051177148
$ws.Cells.Item(8, 20).Value = 0.03903807051177148
$ws.Cells.Item(9, 7).Value = 0.035752
$ws.Cells.Item(9, 8).Value = 0.107256
$ws.Cells.Item(9, 9).Value = 0.2218434834407499
$ws.Cells.Item(9, 10).Value = 0.2218434834407499
$ws.Cells.Item(9, 15).Value = 0.4743638053196239
$ws.Cells.Item(9, 16).Value = 0.4743638053196239
$ws.Cells.Item(9, 17).Value = 0.151509861352
$ws.Cells.Item(9, 18).Value = 1.363588752168
$ws.Cells.Item(9, 19).Value = 0.1052345189903151
$ws.Cells.Item(9, 20).Value = 0.1052345189903151
$ws.Cells.Item(10, 7).Value = 0.035752
$ws.Cells.Item(10, 8).Value = 0.107256
$ws.Cells.Item(10, 9).Value = 0.2218434834407499
$ws.Cells.Item(10, 10).Value = 0.2218434834407499
$ws.Cells.Item(10, 13).Value = 3.123785
$ws.Cells.Item(10, 14).Value = 9.371354999999999
$ws.Cells.Item(10, 15).Value = 0.3496649652969456
$ws.Cells.Item(10, 16).Value = 0.3496649652969455
$ws.Cells.Item(10, 17).Value = 0.11168156132
$ws.Cells.Item(10, 18).Value = 1.00513405188
$ws.Cells.Item(10, 19).Value = 0.07757089393866333
$ws.Cells.Item(10, 20).Value = 0.07757089393866333
